$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update password from Test123 to Test@123 (column B has quotePrefix style - use leading apostrophe to preserve style)
$ws.Cells.Item(2,2).Formula = "'Test@123"

# Row 3: update email/password/exp_result
$ws.Cells.Item(3,1).Formula = "'lgstester@gmail.com"
$ws.Cells.Item(3,2).Formula = "'123233"
$ws.Cells.Item(3,3).Formula = "invalid"

# Delete row 4 entirely
$ws.Range("A4:C4").EntireRow.Delete()
